# Added Mounted Tire Processing Pipeline
# Re-processed tire signal data through the (now mounted-tire-aware) pipeline.
# This refreshes the per-segment probability distributions (Step1_Data),
# their cumulative sums (Step2_Sj), and the derived threshold-crossing
# statistics (Step3_DataPts_*) with the recomputed values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("AN3").Value = 0.003138705231180043
$ws.Range("AO3").Value = 0.08655698538994334
$ws.Range("AP3").Value = 0.1909360282229214
$ws.Range("AQ3").Value = 0.0500555835761112
$ws.Range("AR3").Value = 0.000200943660583779
$ws.Range("AS3").Value = 0.02616634835246307
$ws.Range("AT3").Value = 0.03631088770531193
$ws.Range("AU3").Value = 0.02939179807097527
$ws.Range("AV3").Value = 0.00005183437697472861
$ws.Range("AW3").Value = 0.1429839853810709
$ws.Range("AX3").Value = 0.0002820622239714776
$ws.Range("AY3").Value = 0.01590876521151258
$ws.Range("AZ3").Value = 0.03176929641308554
$ws.Range("BA3").Value = 0.05308779820242941
$ws.Range("BB3").Value = 0.0027036998202216
$ws.Range("BC3").Value = 0.06837097719790038
$ws.Range("BD3").Value = 0.05030200671793431
$ws.Range("BE3").Value = 0.006575261376087838
$ws.Range("BF3").Value = 0.03583498510308193
$ws.Range("BG3").Value = 0.02584112686131336
$ws.Range("BH3").Value = 0.07308735873637592
$ws.Range("BI3").Value = 0.01034137487897129
$ws.Range("BJ3").Value = 0.0007739216739292636
$ws.Range("BK3").Value = 0.004061671156618318
$ws.Range("BL3").Value = 0.01009427021238709
$ws.Range("BM3").Value = 0.006986810648610118
$ws.Range("BN3").Value = 0.0005785058402212682
$ws.Range("BO3").Value = 0.0003475780609990308
$ws.Range("BP3").Value = 0.002874042448471511
$ws.Range("BQ3").Value = 0.01328137786724677
$ws.Range("BR3").Value = 0.01440284535566573
$ws.Range("BS3").Value = 0.006701164025429616
$ws.Range("BT3").Value = 0
$ws.Range("AN5").Value = 0.009561746999716953
$ws.Range("AO5").Value = 0.1413151689822015
$ws.Range("AP5").Value = 0.2185405809949542
$ws.Range("AQ5").Value = 0.1392601343701464
$ws.Range("AR5").Value = 0.0009545019007839453
$ws.Range("AS5").Value = 0.02007906751911928
$ws.Range("AT5").Value = 0.01060844538161553
$ws.Range("AU5").Value = 0.07132452436701017
$ws.Range("AV5").Value = 0.0007861109786253176
$ws.Range("AW5").Value = 0.06073364835199156
$ws.Range("AX5").Value = 0.001651392958055575
$ws.Range("AY5").Value = 0.01983001736255549
$ws.Range("AZ5").Value = 0.0004103044346846094
$ws.Range("BA5").Value = 0.01140721613411309
$ws.Range("BB5").Value = 0.001018945164720214
$ws.Range("BC5").Value = 0.01413965388959351
$ws.Range("BD5").Value = 0.02788096451818821
$ws.Range("BE5").Value = 0.003580921813594172
$ws.Range("BF5").Value = 0.006444133636437301
$ws.Range("BG5").Value = 0.01971034361017913
$ws.Range("BH5").Value = 0.07869422179727012
$ws.Range("BI5").Value = 0.007465030554477349
$ws.Range("BJ5").Value = 0.006064356386951943
$ws.Range("BK5").Value = 0.004550219737737354
$ws.Range("BL5").Value = 0.008171735104729098
$ws.Range("BM5").Value = 0.01547839075234764
$ws.Range("BN5").Value = 0.0002821823637766108
$ws.Range("BO5").Value = 0.0003000737714068087
$ws.Range("BP5").Value = 0.01402062435812443
$ws.Range("BQ5").Value = 0.03788077710682888
$ws.Range("BR5").Value = 0.0406756360269646
$ws.Range("BS5").Value = 0.007178928671098863
$ws.Range("BT5").Value = 0

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("AN3").Value = 0.003138705231180043
$ws.Range("AO3").Value = 0.08969569062112338
$ws.Range("AP3").Value = 0.2806317188440448
$ws.Range("AQ3").Value = 0.330687302420156
$ws.Range("AR3").Value = 0.3308882460807397
$ws.Range("AS3").Value = 0.3570545944332028
$ws.Range("AT3").Value = 0.3933654821385147
$ws.Range("AU3").Value = 0.42275728020949
$ws.Range("AV3").Value = 0.4228091145864647
$ws.Range("AW3").Value = 0.5657930999675356
$ws.Range("AX3").Value = 0.566075162191507
$ws.Range("AY3").Value = 0.5819839274030196
$ws.Range("AZ3").Value = 0.6137532238161052
$ws.Range("BA3").Value = 0.6668410220185346
$ws.Range("BB3").Value = 0.6695447218387561
$ws.Range("BC3").Value = 0.7379156990366565
$ws.Range("BD3").Value = 0.7882177057545908
$ws.Range("BE3").Value = 0.7947929671306787
$ws.Range("BF3").Value = 0.8306279522337606
$ws.Range("BG3").Value = 0.856469079095074
$ws.Range("BH3").Value = 0.9295564378314499
$ws.Range("BI3").Value = 0.9398978127104212
$ws.Range("BJ3").Value = 0.9406717343843505
$ws.Range("BK3").Value = 0.9447334055409689
$ws.Range("BL3").Value = 0.9548276757533559
$ws.Range("BM3").Value = 0.961814486401966
$ws.Range("BN3").Value = 0.9623929922421872
$ws.Range("BO3").Value = 0.9627405703031863
$ws.Range("BP3").Value = 0.9656146127516578
$ws.Range("BQ3").Value = 0.9788959906189046
$ws.Range("BR3").Value = 0.9932988359745704
$ws.Range("BS3").Value = 1
$ws.Range("BT3").Value = 1
$ws.Range("AN5").Value = 0.009561746999716953
$ws.Range("AO5").Value = 0.1508769159819184
$ws.Range("AP5").Value = 0.3694174969768727
$ws.Range("AQ5").Value = 0.508677631347019
$ws.Range("AR5").Value = 0.509632133247803
$ws.Range("AS5").Value = 0.5297112007669222
$ws.Range("AT5").Value = 0.5403196461485378
$ws.Range("AU5").Value = 0.611644170515548
$ws.Range("AV5").Value = 0.6124302814941733
$ws.Range("AW5").Value = 0.6731639298461649
$ws.Range("AX5").Value = 0.6748153228042205
$ws.Range("AY5").Value = 0.694645340166776
$ws.Range("AZ5").Value = 0.6950556446014606
$ws.Range("BA5").Value = 0.7064628607355737
$ws.Range("BB5").Value = 0.707481805900294
$ws.Range("BC5").Value = 0.7216214597898875
$ws.Range("BD5").Value = 0.7495024243080757
$ws.Range("BE5").Value = 0.7530833461216698
$ws.Range("BF5").Value = 0.7595274797581071
$ws.Range("BG5").Value = 0.7792378233682862
$ws.Range("BH5").Value = 0.8579320451655563
$ws.Range("BI5").Value = 0.8653970757200337
$ws.Range("BJ5").Value = 0.8714614321069856
$ws.Range("BK5").Value = 0.876011651844723
$ws.Range("BL5").Value = 0.8841833869494521
$ws.Range("BM5").Value = 0.8996617777017998
$ws.Range("BN5").Value = 0.8999439600655764
$ws.Range("BO5").Value = 0.9002440338369833
$ws.Range("BP5").Value = 0.9142646581951077
$ws.Range("BQ5").Value = 0.9521454353019366
$ws.Range("BR5").Value = 0.9928210713289012
$ws.Range("BS5").Value = 1

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F3").Value = 0.5657930999675356
$ws.Range("D5").Value = 42
$ws.Range("F5").Value = 0.508677631347019
$ws.Range("G5").Value = 4

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F3").Value = 0.7379156990366565
$ws.Range("D5").Value = 52
$ws.Range("F5").Value = 0.7064628607355737
$ws.Range("G5").Value = 14

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F3").Value = 0.8306279522337606
$ws.Range("F5").Value = 0.8579320451655563

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F3").Value = 0.9295564378314499
$ws.Range("D5").Value = 66
$ws.Range("F5").Value = 0.9002440338369833
$ws.Range("G5").Value = 28
